$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2042.3077
$ws.Range("I40").Value = 1894.4445
$ws.Range("K40").Value = 1894.4445
$ws.Range("M40").Value = -1719.4445
$ws.Range("H125").Value = 3153
$ws.Range("J125").Value = 3694.4
$ws.Range("L125").Value = 33249.6
$ws.Range("N125").Value = -38169.6
$ws.Range("H129").Value = 3893.4614
$ws.Range("I129").Value = 2033
$ws.Range("K129").Value = 6099
$ws.Range("M129").Value = -1099
$ws.Range("H137").Value = 2781.1052
$ws.Range("J137").Value = 5358.5713
$ws.Range("L137").Value = 16075.7139
$ws.Range("N137").Value = -21175.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12813.75
$ws.Range("I32").Value = 8751.75
$ws.Range("K32").Value = 8751.75
$ws.Range("M32").Value = -8464.75
$ws.Range("H97").Value = 729.5294
$ws.Range("I97").Value = 482.66666
$ws.Range("J97").Value = 2581
$ws.Range("K97").Value = 482.66666
$ws.Range("L97").Value = 2581
$ws.Range("M97").Value = 13.33334000000002
$ws.Range("N97").Value = -3573
$ws.Range("H122").Value = 2735.35
$ws.Range("I122").Value = 1683.9166
$ws.Range("J122").Value = 4312.5
$ws.Range("K122").Value = 5051.7498
$ws.Range("L122").Value = 12937.5
$ws.Range("M122").Value = -2601.7498
$ws.Range("N122").Value = -17837.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9288.6
$ws.Range("I20").Value = 9654.111000000001
$ws.Range("K20").Value = 9654.111000000001
$ws.Range("M20").Value = -9407.111000000001
$ws.Range("H22").Value = 670.2
$ws.Range("I22").Value = 670.2
$ws.Range("K22").Value = 670.2
$ws.Range("M22").Value = -497.2
$ws.Range("H99").Value = 2216.5945
$ws.Range("J99").Value = 3260.3
$ws.Range("L99").Value = 3260.3
$ws.Range("N99").Value = -6256.3
$ws.Range("H103").Value = 7582.2
$ws.Range("J103").Value = 7582.2
$ws.Range("L103").Value = 7582.2
$ws.Range("N103").Value = -9926.200000000001
$ws.Range("H134").Value = 2071.8845
$ws.Range("I134").Value = 1794.1364
$ws.Range("K134").Value = 5382.4092
$ws.Range("M134").Value = -2847.4092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 170.52174
$ws.Range("I7").Value = 139.19048
$ws.Range("K7").Value = 139.19048
$ws.Range("M7").Value = -26.19048000000001
$ws.Range("H16").Value = 685.2941
$ws.Range("I16").Value = 680.9167
$ws.Range("K16").Value = 680.9167
$ws.Range("M16").Value = -393.9167
$ws.Range("H74").Value = 46109.4
$ws.Range("J74").Value = 45074.25
$ws.Range("L74").Value = 45074.25
$ws.Range("N74").Value = -46822.25
$ws.Range("H77").Value = 46109.4
$ws.Range("J77").Value = 45074.25
$ws.Range("L77").Value = 135222.75
$ws.Range("N77").Value = -143958.75
$ws.Range("H99").Value = 14017.654
$ws.Range("I99").Value = 11096.272
$ws.Range("K99").Value = 11096.272
$ws.Range("M99").Value = -9598.272000000001
$ws.Range("H102").Value = 16075
$ws.Range("I102").Value = 12150
$ws.Range("K102").Value = 12150
$ws.Range("M102").Value = -9716
$ws.Range("H104").Value = 9000
$ws.Range("J104").Value = 9000
$ws.Range("L104").Value = 9000
$ws.Range("N104").Value = -14242
$ws.Range("H113").Value = 685.2941
$ws.Range("I113").Value = 680.9167
$ws.Range("K113").Value = 680.9167
$ws.Range("M113").Value = 1489.0833
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 2132.2693
$ws.Range("I122").Value = 2190.2173
$ws.Range("K122").Value = 6570.651899999999
$ws.Range("M122").Value = -4120.651899999999
$ws.Range("H126").Value = 14017.654
$ws.Range("I126").Value = 11096.272
$ws.Range("K126").Value = 33288.81600000001
$ws.Range("M126").Value = -30818.81600000001
$ws.Range("H134").Value = 2042.5143
$ws.Range("I134").Value = 1541.32
$ws.Range("K134").Value = 4623.96
$ws.Range("M134").Value = -2088.96

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 523.625
$ws.Range("J5").Value = 552.5
$ws.Range("L5").Value = 1657.5
$ws.Range("N5").Value = -1881.5
$ws.Range("H8").Value = 186.83333
$ws.Range("I8").Value = 186.83333
$ws.Range("K8").Value = 560.49999
$ws.Range("M8").Value = -421.49999
$ws.Range("H113").Value = 1172.875
$ws.Range("J113").Value = 1638.6
$ws.Range("L113").Value = 4915.799999999999
$ws.Range("N113").Value = -9255.799999999999
$ws.Range("H119").Value = 2781
$ws.Range("I119").Value = 2781
$ws.Range("K119").Value = 8343
$ws.Range("M119").Value = -3505
$ws.Range("H135").Value = 523.625
$ws.Range("J135").Value = 552.5
$ws.Range("L135").Value = 4972.5
$ws.Range("N135").Value = -10042.5
$ws.Range("H137").Value = 4650.8667
$ws.Range("J137").Value = 7290.8
$ws.Range("L137").Value = 21872.4
$ws.Range("N137").Value = -32072.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7287.25
$ws.Range("I70").Value = 5808
$ws.Range("K70").Value = 5808
$ws.Range("M70").Value = -5538
$ws.Range("H73").Value = 7287.25
$ws.Range("I73").Value = 5808
$ws.Range("K73").Value = 5808
$ws.Range("M73").Value = -4872
$ws.Range("H98").Value = 17229
$ws.Range("J98").Value = 17229
$ws.Range("L98").Value = 17229
$ws.Range("N98").Value = -23219
$ws.Range("H102").Value = 2179.8147
$ws.Range("J102").Value = 3857.0715
$ws.Range("L102").Value = 3857.0715
$ws.Range("N102").Value = -7101.0715
$ws.Range("H113").Value = 3665.7144
$ws.Range("I113").Value = 4137
$ws.Range("J113").Value = 3518.4375
$ws.Range("K113").Value = 4137
$ws.Range("L113").Value = 3518.4375
$ws.Range("M113").Value = -1967
$ws.Range("N113").Value = -7858.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2578.7273
$ws.Range("I7").Value = 2596.2222
$ws.Range("K7").Value = 2596.2222
$ws.Range("M7").Value = -2484.2222
$ws.Range("H22").Value = 8514
$ws.Range("J22").Value = 17833
$ws.Range("L22").Value = 17833
$ws.Range("N22").Value = -18423
$ws.Range("H27").Value = 8514
$ws.Range("J27").Value = 17833
$ws.Range("L27").Value = 17833
$ws.Range("N27").Value = -18047
$ws.Range("H40").Value = 3100.6667
$ws.Range("I40").Value = 2843.8572
$ws.Range("K40").Value = 2843.8572
$ws.Range("M40").Value = -2707.8572
$ws.Range("H46").Value = 3409.25
$ws.Range("I46").Value = 2456.8572
$ws.Range("J46").Value = 3922.077
$ws.Range("K46").Value = 2456.8572
$ws.Range("L46").Value = 3922.077
$ws.Range("M46").Value = -2268.8572
$ws.Range("N46").Value = -4298.077
$ws.Range("H122").Value = 8839
$ws.Range("I122").Value = 8932.223
$ws.Range("K122").Value = 26796.669
$ws.Range("M122").Value = -24346.669
$ws.Range("H124").Value = 62499.5
$ws.Range("J124").Value = 62499.5
$ws.Range("L124").Value = 62499.5
$ws.Range("N124").Value = -72319.5
$ws.Range("H126").Value = 2578.7273
$ws.Range("I126").Value = 2596.2222
$ws.Range("K126").Value = 7788.6666
$ws.Range("M126").Value = -5318.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 271996.75
$ws.Range("I3").Value = 352995.66
$ws.Range("K3").Value = 352995.66
$ws.Range("M3").Value = -352881.66
$ws.Range("H4").Value = 599.8
$ws.Range("I4").Value = 999
$ws.Range("K4").Value = 999
$ws.Range("M4").Value = -886
$ws.Range("H126").Value = 2337.2
$ws.Range("J126").Value = 4001.6667
$ws.Range("K126").Value = 4001.6667
$ws.Range("L126").Value = 12005.0001
$ws.Range("N126").Value = -16945.0001

Write-Output "Applied changes to $([int]201) cells across 8 sheets"